$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Row 16: existing row gets the B (status) formula plus new K/L/M cells
# ---------------------------------------------------------------------
$ws.Range("B16").Formula = "=IF(E16=`"`",0,IF(F16=`"`",1,IF(H16=`"`",2,3)))"
$ws.Range("K16").Value = 1

# ---------------------------------------------------------------------
# Row 17: SATA Data Connector
# ---------------------------------------------------------------------
$ws.Range("B17").Formula = "=IF(E17=`"`",0,IF(F17=`"`",1,IF(H17=`"`",2,3)))"
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = "SATA Data Connector"
$ws.Range("E17").Value = "Connector for Signals of SATA drive"
$ws.Range("F17").Value = "Molex"
$ws.Range("G17").Value = 471554001
$ws.Range("G17").HorizontalAlignment = -4131
$ws.Range("K17").Value = 5

# ---------------------------------------------------------------------
# Row 18: Gigabit Socket
# ---------------------------------------------------------------------
$ws.Range("B18").Formula = "=IF(E18=`"`",0,IF(F18=`"`",1,IF(H18=`"`",2,3)))"
$ws.Range("C18").Value = 13
$ws.Range("D18").Value = "Gigabit Socket"
$ws.Range("K18").Value = 1

# ---------------------------------------------------------------------
# Row 19: USB Connector
# ---------------------------------------------------------------------
$ws.Range("B19").Formula = "=IF(E19=`"`",0,IF(F19=`"`",1,IF(H19=`"`",2,3)))"
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = "USB Connector"
$ws.Range("K19").Value = 12

# ---------------------------------------------------------------------
# Row 20: BIOS EEPROM
# ---------------------------------------------------------------------
$ws.Range("B20").Formula = "=IF(E20=`"`",0,IF(F20=`"`",1,IF(H20=`"`",2,3)))"
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = "BIOS EEPROM"
$ws.Range("E20").Value = "16Mbit W25X16BV EEPROM"
$ws.Range("F20").Value = "Winbond"
$ws.Range("G20").Value = "W25X16BVDAI"
$ws.Range("K20").Value = 1

# ---------------------------------------------------------------------
# L / M columns: fill L16:L19 / M16:M19 together so they form one shared
# formula group (matches how the extension was authored), then L20/M20
# separately.
# ---------------------------------------------------------------------
$ws.Range("L16:L19").Formula = "=K16*`$K`$4"
$ws.Range("M16:M19").Formula = "=L16*J16"
$ws.Range("L20").Formula = "=K20*`$K`$4"
$ws.Range("M20").Formula = "=L20*J20"

# ---------------------------------------------------------------------
# Conditional formatting: extend both rules to cover the new rows.
# ---------------------------------------------------------------------

# Rule 1 (expression over "A7:B7 B8:B15") -> "A7:B7 B8:B20"
$exprRule = $null
$fcsB7 = $ws.Range("B7").FormatConditions
for ($i = 1; $i -le $fcsB7.Count; $i++) {
    $r = $fcsB7.Item($i)
    if ($r.Type -eq 2) { $exprRule = $r }
}
$exprFormula = $exprRule.Formula1
$exprRule.Delete()
$newExprA = $ws.Range("A7:B7").FormatConditions.Add(2, 0, $exprFormula)
$newExprA.Priority = 2
$newExprB = $ws.Range("B8:B20").FormatConditions.Add(2, 0, $exprFormula)
$newExprB.Priority = 2

# Rule 2 (iconSet over "B7:B15") -> "B7:B20" with priority 6
$fcsB7b = $ws.Range("B7").FormatConditions
for ($i = 1; $i -le $fcsB7b.Count; $i++) {
    $r = $fcsB7b.Item($i)
    if ($r.Type -eq 6) {
        $r.ModifyAppliesToRange($ws.Range("B7:B20"))
        $r.Priority = 6
    }
}

# ---------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------
$ws.Range("H18").Select()
